# Applies the weekly refresh of Fruta/Chirimoya pricing data.
# The underlying rows were re-sorted; this script writes the resulting
# cell values directly onto the existing rows (2-17) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44848
$ws.Range("L2").Value = "Especial"
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 24500
$ws.Range("S2").Value = 2450

# Row 3
$ws.Range("D3").Value = 44848
$ws.Range("M3").Value = 120

# Row 4
$ws.Range("D4").Value = 44461
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 31000
$ws.Range("O4").Value = 32000
$ws.Range("P4").Value = 31500
$ws.Range("S4").Value = 3150

# Row 5
$ws.Range("D5").Value = 44461
$ws.Range("L5").Value = "Primera"
$ws.Range("N5").Value = 30000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 30000
$ws.Range("S5").Value = 3000

# Row 6
$ws.Range("D6").Value = 44841
$ws.Range("M6").Value = 60

# Row 7
$ws.Range("D7").Value = 44874
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 25000
$ws.Range("S7").Value = 2500

# Row 8
$ws.Range("D8").Value = 44874
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 24000
$ws.Range("P8").Value = 23500
$ws.Range("S8").Value = 2350

# Row 10
$ws.Range("D10").Value = 44487
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("S10").Value = 2350

# Row 12
$ws.Range("D12").Value = 44448

# Row 15
$ws.Range("D15").Value = 44452

# Row 16
$ws.Range("D16").Value = 44839
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 120
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 26000
$ws.Range("P16").Value = 25500
$ws.Range("S16").Value = 2550

# Row 17
$ws.Range("D17").Value = 44447
$ws.Range("M17").Value = 60
